# ---------------------------------------------------------------------------
# Edits the "journal of incoming/outgoing documents" table:
#   1. Clears the "№14351дск" run from row 1 (col 2), leaving an empty
#      centered paragraph.
#   2. Renames the document title in row 1 (col 4) to "Ключ".
#   3. Renames the document title in row 2 (col 4) to "зявление".
#   4. Appends nine new data rows (index 3, blank, 5, blank, 7, 8, 9, 10, 11)
#      describing further documents / sheet ranges.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

$quoteOpen  = [char]8220   # “
$quoteClose = [char]8221   # ”

# --- 1. remove the "№14351дск" run entirely --------------------------------
$d.Content.Find.Execute("№14351дск", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# --- 2. "Заявка директора завода "Древострой"" -> "Ключ" -------------------
$oldTitle1 = "Заявка директора завода " + $quoteOpen + "Древострой" + $quoteClose
$d.Content.Find.Execute($oldTitle1, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ключ", 2) | Out-Null

# --- 3. "на древесину, от 02.02.2022" -> "зявление" -------------------------
$d.Content.Find.Execute("на древесину, от 02.02.2022", $true, $false, $false, `
                         $false, $false, $true, 1, $false, "зявление", 2) | Out-Null

# --- 4. append the new rows --------------------------------------------------
$t = $d.Tables.Item(1)

# Always sets the alignment explicitly (1 = center, 0 = left/default) so that
# formatting never leaks in from whichever row Word used as the template for
# the freshly-added row.
function Fill-Cell($table, $rowIndex, $colIndex, $text, $center) {
    $cell = $table.Cell($rowIndex, $colIndex)
    if ($text) {
        $cell.Range.Text = $text
    }
    if ($center) {
        $cell.Range.ParagraphFormat.Alignment = 1
    } else {
        $cell.Range.ParagraphFormat.Alignment = 0
    }
}

function Add-DataRow($table, $num, $col4parts, $col5) {
    $row = $table.Rows.Add()
    $idx = $row.Index

    $numCentered = ($num -ne $null -and $num -ne "")

    # column 1: the row number (centered when present, plain otherwise)
    Fill-Cell $table $idx 1 $num $numCentered

    # column 2: always empty, but keeps the same centering as column 1
    Fill-Cell $table $idx 2 $null $numCentered

    if ($col4parts.Count -gt 0) {
        $cell4 = $table.Cell($idx, 4)
        $cell4.Range.Text = $col4parts[0]
        for ($i = 1; $i -lt $col4parts.Count; $i++) {
            $endRange = $cell4.Range
            $endRange.Collapse(0)
            $endRange.InsertAfter($col4parts[$i])
        }
    }

    # column 5: sheet-number range (centered when present, plain otherwise)
    $col5Centered = ($col5 -ne $null -and $col5 -ne "")
    Fill-Cell $table $idx 5 $col5 $col5Centered

    return $idx
}

# Row: "3" | Ключ / Одна строка схема, от 21.04.2006 | 2-4
Add-DataRow $t "3" @("Ключ", "Одна строка схема, от 21.04.2006") "2-4" | Out-Null

# Row: (blank) | зявление, от 21.04.2006
Add-DataRow $t $null @("зявление, от 21.04.2006") $null | Out-Null

# Row: "5" | Ключ | 5-7
Add-DataRow $t "5" @("Ключ") "5-7" | Out-Null

# Row: (blank) | зявление, от 21.05.2006
Add-DataRow $t $null @("зявление, от 21.05.2006") $null | Out-Null

# Row: "7" | Одна строка схема | 8-10
Add-DataRow $t "7" @("Одна строка схема") "8-10" | Out-Null

# Row: "8" | Одна строка схема, от 21.07.2006 | 11-15
Add-DataRow $t "8" @("Одна строка схема, от 21.07.2006") "11-15" | Out-Null

# Row: "9" | (blank) | 16-25
Add-DataRow $t "9" @() "16-25" | Out-Null

# Row: "10" | (blank) | 26-28
Add-DataRow $t "10" @() "26-28" | Out-Null

# Row: "11" | (blank) | 29-30
Add-DataRow $t "11" @() "29-30" | Out-Null

Write-Host "Edit complete. Table now has" $t.Rows.Count "rows."
